$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 16:33"

# Update the Gran Canaria row (row 31) figures
$ws.Range("B31").Value = 2235
$ws.Range("C31").Value = 1232
$ws.Range("D31").Value = 859
$ws.Range("E31").Value = 144
